$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------------
# Overview sheet: status columns (zh-cn / de-de) for both rows move from
# "Ready for handoff" to the handed-back status text.
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = $statusText
$wsOverview.Range("C2").Value = $statusText
$wsOverview.Range("B3").Value = $statusText
$wsOverview.Range("C3").Value = $statusText

# ---------------------------------------------------------------------------
# zh-cn sheet: mark status as handed back, and fill in the "Latest Target
# File" / "Latest Handback File" columns (F/G) for both data rows.
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = $statusText
$wsZh.Range("C3").Value = $statusText

$wsZh.Range("F2").Value = "1cdee9db-87a8-432b-93c7-7a1de8ca5c9d.md"
$wsZh.Hyperlinks.Add($wsZh.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/96d8c3d4c9fd08ac52b4e271e62feacf2fb416e2/e2e/1cdee9db-87a8-432b-93c7-7a1de8ca5c9d.md", "", "1cdee9db-87a8-432b-93c7-7a1de8ca5c9d.md", "1cdee9db-87a8-432b-93c7-7a1de8ca5c9d.md") | Out-Null

$wsZh.Range("G2").Value = "1cdee9db-87a8-432b-93c7-7a1de8ca5c9d.e32f7de97cb30ca65613d8b7908b42c19a91097e.zh-cn.xlf"
$wsZh.Hyperlinks.Add($wsZh.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/dc4f301059e22ef9cd0842235d794284d783e668/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/1cdee9db-87a8-432b-93c7-7a1de8ca5c9d.e32f7de97cb30ca65613d8b7908b42c19a91097e.zh-cn.xlf", "", "1cdee9db-87a8-432b-93c7-7a1de8ca5c9d.e32f7de97cb30ca65613d8b7908b42c19a91097e.zh-cn.xlf", "1cdee9db-87a8-432b-93c7-7a1de8ca5c9d.e32f7de97cb30ca65613d8b7908b42c19a91097e.zh-cn.xlf") | Out-Null

$wsZh.Range("F3").Value = "ffa9b05e-361c-4efb-966a-0babefbfbb56.md"
$wsZh.Hyperlinks.Add($wsZh.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/96d8c3d4c9fd08ac52b4e271e62feacf2fb416e2/e2e/ffa9b05e-361c-4efb-966a-0babefbfbb56.md", "", "ffa9b05e-361c-4efb-966a-0babefbfbb56.md", "ffa9b05e-361c-4efb-966a-0babefbfbb56.md") | Out-Null

$wsZh.Range("G3").Value = "ffa9b05e-361c-4efb-966a-0babefbfbb56.10bd8a104e27074121eb03c63a1d3a659f523cd6.zh-cn.xlf"
$wsZh.Hyperlinks.Add($wsZh.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/dc4f301059e22ef9cd0842235d794284d783e668/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/ffa9b05e-361c-4efb-966a-0babefbfbb56.10bd8a104e27074121eb03c63a1d3a659f523cd6.zh-cn.xlf", "", "ffa9b05e-361c-4efb-966a-0babefbfbb56.10bd8a104e27074121eb03c63a1d3a659f523cd6.zh-cn.xlf", "ffa9b05e-361c-4efb-966a-0babefbfbb56.10bd8a104e27074121eb03c63a1d3a659f523cd6.zh-cn.xlf") | Out-Null

# ---------------------------------------------------------------------------
# de-de sheet: same as zh-cn, plus this locale's handback finished, so the
# "Latest Handback DateTime" column (H) gets the real timestamp.
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = $statusText
$wsDe.Range("C3").Value = $statusText

$wsDe.Range("F2").Value = "1cdee9db-87a8-432b-93c7-7a1de8ca5c9d.md"
$wsDe.Hyperlinks.Add($wsDe.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/96d8c3d4c9fd08ac52b4e271e62feacf2fb416e2/e2e/1cdee9db-87a8-432b-93c7-7a1de8ca5c9d.md", "", "1cdee9db-87a8-432b-93c7-7a1de8ca5c9d.md", "1cdee9db-87a8-432b-93c7-7a1de8ca5c9d.md") | Out-Null

$wsDe.Range("G2").Value = "1cdee9db-87a8-432b-93c7-7a1de8ca5c9d.e32f7de97cb30ca65613d8b7908b42c19a91097e.de-de.xlf"
$wsDe.Hyperlinks.Add($wsDe.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/3e75ae2b235c95bb810ba38271a308836752ff5c/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/1cdee9db-87a8-432b-93c7-7a1de8ca5c9d.e32f7de97cb30ca65613d8b7908b42c19a91097e.de-de.xlf", "", "1cdee9db-87a8-432b-93c7-7a1de8ca5c9d.e32f7de97cb30ca65613d8b7908b42c19a91097e.de-de.xlf", "1cdee9db-87a8-432b-93c7-7a1de8ca5c9d.e32f7de97cb30ca65613d8b7908b42c19a91097e.de-de.xlf") | Out-Null

$wsDe.Range("F3").Value = "ffa9b05e-361c-4efb-966a-0babefbfbb56.md"
$wsDe.Hyperlinks.Add($wsDe.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/96d8c3d4c9fd08ac52b4e271e62feacf2fb416e2/e2e/ffa9b05e-361c-4efb-966a-0babefbfbb56.md", "", "ffa9b05e-361c-4efb-966a-0babefbfbb56.md", "ffa9b05e-361c-4efb-966a-0babefbfbb56.md") | Out-Null

$wsDe.Range("G3").Value = "ffa9b05e-361c-4efb-966a-0babefbfbb56.10bd8a104e27074121eb03c63a1d3a659f523cd6.de-de.xlf"
$wsDe.Hyperlinks.Add($wsDe.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/3e75ae2b235c95bb810ba38271a308836752ff5c/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/ffa9b05e-361c-4efb-966a-0babefbfbb56.10bd8a104e27074121eb03c63a1d3a659f523cd6.de-de.xlf", "", "ffa9b05e-361c-4efb-966a-0babefbfbb56.10bd8a104e27074121eb03c63a1d3a659f523cd6.de-de.xlf", "ffa9b05e-361c-4efb-966a-0babefbfbb56.10bd8a104e27074121eb03c63a1d3a659f523cd6.de-de.xlf") | Out-Null

$wsDe.Range("H2").Value = "2016-03-11 14:32:35"
$wsDe.Range("H3").Value = "2016-03-11 14:32:35"

Write-Host "Handback report generated."
